$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.132.07'
$ws.Range("E2").Value = '  +0.45%  '

$ws.Range("D3").Value = '2.354.86'
$ws.Range("E3").Value = '  -0.16%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("E5").Value = '  +1.89%  '

$ws.Range("D6").Value = "'240.03"
$ws.Range("E6").Value = '  +1.83%  '

$ws.Range("D7").Value = "'74.73"
$ws.Range("E7").Value = '  +1.63%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").Value = "'0.602"
$ws.Range("E9").Value = '  +6.99%  '

$ws.Range("E10").Value = '  +2.45%  '

$ws.Range("D11").Value = "'57.18"
$ws.Range("E11").Value = '  +0.15%  '

$ws.Range("D12").Value = "'32.32"
$ws.Range("E12").Value = '  +15.15%  '

$ws.Range("E13").Value = '  +7.59%  '

$ws.Range("E14").Value = '  +0.96%  '

$ws.Range("D15").Value = '2.693.10'
$ws.Range("E15").Value = '  -0.54%  '

$ws.Range("D16").Value = "'16.56"
$ws.Range("E16").Value = '  -1.64%  '

$ws.Range("D17").Value = "'0.909"
$ws.Range("E17").Value = '  +2.20%  '

$ws.Range("D18").Value = '2.354.50'
$ws.Range("E18").Value = '  -1.70%  '

$ws.Range("D19").Value = '43.940.70'
$ws.Range("E19").Value = '  +0.00%  '

$ws.Range("E20").Value = '  +1.53%  '

$ws.Range("D21").Value = "'6.73"
$ws.Range("E21").Value = '  +4.94%  '

$ws.Range("D22").Value = "'77.15"
$ws.Range("E22").Value = '  -1.39%  '

$ws.Range("D23").Value = "'257.23"
$ws.Range("E23").Value = '  +1.44%  '

$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("D25").Value = "'1.88"
$ws.Range("E25").Value = '  +18.57%  '

$ws.Range("E26").Value = '  -1.36%  '

$ws.Range("D27").Value = "'2.50"
$ws.Range("E27").Value = '  +0.62%  '

$ws.Range("D28").Value = "'10.75"
$ws.Range("E28").Value = '  +0.57%  '

$ws.Range("E29").Value = '  -1.72%  '

$ws.Range("D30").Value = "'22.91"
$ws.Range("E30").Value = '  +1.84%  '

$ws.Range("D31").Value = "'174.93"
$ws.Range("E31").Value = '  +1.35%  '

$ws.Range("E32").Value = '  -2.16%  '

$ws.Range("D33").Value = "'0.138"
$ws.Range("E33").Value = '  +3.60%  '

$ws.Range("D34").Value = "'0.0761"
$ws.Range("E34").Value = '  +5.56%  '

$ws.Range("D35").Value = "'5.32"
$ws.Range("E35").Value = '  +2.24%  '

$ws.Range("E36").Value = '  +3.36%  '

$ws.Range("D37").Value = "'3.75"
$ws.Range("E37").Value = '  -0.51%  '

$ws.Range("E38").Value = '  -2.91%  '

$ws.Range("E39").Value = '  -0.91%  '

$ws.Range("E40").Value = '  +4.00%  '

$ws.Range("D41").Value = "'0.208"
$ws.Range("E41").Value = '  +11.72%  '

$ws.Range("B42").Value = 'Cronos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D42").Value = "'0.109"
$ws.Range("E42").Value = '  +11.11%  '

$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = "'19.41"
$ws.Range("E43").Value = '  +0.90%  '

$ws.Range("D44").Value = "'9.08"
$ws.Range("E44").Value = '  +1.10%  '

$ws.Range("E45").Value = '  +0.08%  '

$ws.Range("D46").Value = "'4.80"
$ws.Range("E46").Value = '  +7.92%  '

$ws.Range("D47").Value = "'2.52"
$ws.Range("E47").Value = '  +9.31%  '

$ws.Range("E48").Value = '  +2.75%  '

$ws.Range("D49").Value = "'101.25"
$ws.Range("E49").Value = '  +3.14%  '

$ws.Range("E50").Value = '  +0.94%  '

$ws.Range("D51").Value = "'56.75"
$ws.Range("E51").Value = '  +9.28%  '
